$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Sort an array of 0s, 1s and 2s"
$ws.Range("B9").Value = "DutchNationalFlag"

$ws.Range("B9").Select()
